$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rubricas")

# Fill in previously empty "D" values for the two newly-scored rows
$ws.Range("D19").Value = 1
$ws.Range("D20").Value = 1

# Update the view so it matches the scrolled/selected position from the edit
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D31").Select()
